$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to make changes, then re-protect afterwards.
$ws.Unprotect()

# Update the confidential-notice text (shared string): date 2021-03-29 -> 2021-03-30.
# Setting a multi-line .Value directly nudges Excel's auto row-height for that row,
# so we restore it with AutoFit to keep the row's height attribute untouched.
$ws.Range("A40").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-30 for illustrative purposes only and are subject to change."
$ws.Rows(40).AutoFit()

# Update the Weight (D) and Percent Change (E) columns for each holding row.
$ws.Range("D2").Value = 0.03246642236391084
$ws.Range("E2").Value = 0.0005621662138104444
$ws.Range("D3").Value = 0.03426322423422986
$ws.Range("E3").Value = -0.008224665473451021
$ws.Range("D4").Value = 0.03191291442674817
$ws.Range("E4").Value = -0.0147783251231528
$ws.Range("D5").Value = 0.03742451016859421
$ws.Range("E5").Value = -0.006645576822412957
$ws.Range("D6").Value = 0.01690212960019059
$ws.Range("E6").Value = -0.004319343459794123
$ws.Range("D7").Value = 0.01597251855602203
$ws.Range("E7").Value = -0.008836748685914442
$ws.Range("D8").Value = 0.03412499934598699
$ws.Range("E8").Value = -0.0204345779730154
$ws.Range("D9").Value = 0.0325373599606059
$ws.Range("E9").Value = -0.0004038772213247332
$ws.Range("D10").Value = 0.03126790550722754
$ws.Range("E10").Value = -0.001089600155657289
$ws.Range("D11").Value = 0.0310331908862588
$ws.Range("E11").Value = -0.01125679289226256
$ws.Range("D12").Value = 0.01495164988317198
$ws.Range("E12").Value = 0.02994791666666674
$ws.Range("D13").Value = 0.01670318796966499
$ws.Range("E13").Value = -0.003642323802586001
$ws.Range("D14").Value = 0.008013149859262485
$ws.Range("E14").Value = 0.01254251700680298
$ws.Range("D15").Value = 0.007845722529841549
$ws.Range("E15").Value = 0.0188585607940448
$ws.Range("D16").Value = 0.03077499263550937
$ws.Range("E16").Value = -0.003653271338425856
$ws.Range("D17").Value = 0.03209786322087595
$ws.Range("E17").Value = -0.003241140881590443
$ws.Range("D18").Value = 0.03261669325913966
$ws.Range("E18").Value = -0.01862269641125103
$ws.Range("D19").Value = 0.031847452287774
$ws.Range("E19").Value = -0.009696719620383765
$ws.Range("D20").Value = 0.02501614955835566
$ws.Range("E20").Value = -0.0001994211921495515
$ws.Range("D21").Value = 0.03170703721644277
$ws.Range("E21").Value = 0.01927977158996708
$ws.Range("D22").Value = 0.03309366646504814
$ws.Range("E22").Value = 0.008669755129053769
$ws.Range("D23").Value = 0.0325782433782552
$ws.Range("E23").Value = 0.01195171507111281
$ws.Range("D24").Value = 0.01613337533617789
$ws.Range("E24").Value = 0.0299263907324725
$ws.Range("D25").Value = 0.015027089522882
$ws.Range("E25").Value = 0.01578947368421058
$ws.Range("D26").Value = 0.03306203048710524
$ws.Range("E26").Value = -0.008390990725747116
$ws.Range("D27").Value = 0.03338812441359369
$ws.Range("E27").Value = 0.002704081632653121
$ws.Range("D28").Value = 0.03148558537084241
$ws.Range("E28").Value = -0.01441081448733228
$ws.Range("D29").Value = 0.03248528227383835
$ws.Range("E29").Value = -0.004045246834968985
$ws.Range("D30").Value = 0.03151004241532904
$ws.Range("E30").Value = -0.005908134303863388
$ws.Range("D31").Value = 0.0334780435970545
$ws.Range("E31").Value = -0.01567571300324566
$ws.Range("D32").Value = 0.0334046724635946
$ws.Range("E32").Value = 0.005944575574788002
$ws.Range("D33").Value = 0.03228597561279802
$ws.Range("E33").Value = -0.01450957632037142
$ws.Range("D34").Value = 0.04724224921582319
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 0.03158195342673004
$ws.Range("E35").Value = 0.01178936337437775
$ws.Range("D36").Value = 0.03376459255111431
$ws.Range("E36").Value = 0.007265021946420314
$ws.Range("E37").Value = -0.001922007336869869

# Re-protect the sheet (best effort - the original protection flags/password hash
# cannot be reproduced exactly through this object model, but we restore the
# "sheet is protected" state to match the source workbook).
$ws.Protect()
